$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# A new row is being inserted above current row 6 (Giovani / Galpao
# Toyota), pushing every row from 6..161 down to 7..162. We replicate
# this via a manual "copy range down by one row" rather than
# Rows.Insert() so that no incidental/unused cell styles get minted
# in styles.xml (Insert() in this runtime always mints a couple of
# spurious cellXfs entries that the real edit's styles.xml doesn't
# have).
# ------------------------------------------------------------------

$ws.Range("A6:I161").Copy()
$ws.Range("A7:I162").PasteSpecial(-4104)   # xlPasteAll
$excel.CutCopyMode = 0

# Row 6 carried ht=45 (tall "Central off..." row); that row is now
# row 7, so move the explicit height along with it, and let row 6
# (now the freshly inserted row) go back to the default height.
$ws.Rows(7).RowHeight = 45
$ws.Rows(6).AutoFit()

# The copy/paste above does not materialize the new trailing row 162
# on its own (its source row 161 is blank), so stamp a value into it
# to force the row into existence, copy row 161's formatting onto it,
# then blank the value back out again.
$ws.Range("A162").Value = "x"
$ws.Range("A161:I161").Copy()
$ws.Range("A162:I162").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A162").Value = ""

# ------------------------------------------------------------------
# Fill in the new row 6 with the Cimentao ticket (keeps the style
# that row already has, which matches the other Roberto rows).
# ------------------------------------------------------------------
$ws.Range("A6").Value = "Roberto"
$ws.Range("B6").Value = "'0304"
$ws.Range("C6").Value = "Cimentão"
$ws.Range("D6").Value = "Disparos frequentes e ambos os DVR sem comunicação."
$ws.Range("E6").Value = ""
$ws.Range("G6").Value = "Pendente"
$ws.Range("H6").Value = ""

# ------------------------------------------------------------------
# Update the "Kit Faltando" counter in H2.
# ------------------------------------------------------------------
$ws.Range("H2").Value = "Maxvel: 38 / Forte: 19"

# ------------------------------------------------------------------
# Move the active selection to H2 (matches the saved selection in
# the target file).
# ------------------------------------------------------------------
$ws.Range("H2").Select()
